# Mark progress for DSA problems #11-#15 (rows 12-16) as Done, filling in
# the Date Done / Language / Time Taken / Difficulty / Notes columns,
# mirroring the pattern already used for the earlier, already-completed rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date Done for all five rows is 2026-02-26 (serial 46079), same as row 11's date style.
$rowsData = @(
    @{ Row = 12; Status = "Done"; DateDone = 46079; Language = "Python"; Time = "2 mins"; Difficulty = "Easy"; Notes = "Used πr², float input, exponent operator" },
    @{ Row = 13; Status = "Done"; DateDone = 46079; Language = "Python"; Time = "3 mins"; Difficulty = "Easy"; Notes = "Used lower(), isalpha(), membership operator(in)" },
    @{ Row = 14; Status = "Done"; DateDone = 46079; Language = "Python"; Time = "2 mins"; Difficulty = "Easy"; Notes = "Used ord() to convert a character to an ASCII value; validated single-character input" },
    @{ Row = 15; Status = "Done"; DateDone = 46079; Language = "Python"; Time = "2 mins"; Difficulty = "Easy"; Notes = "Used modulo operator; used logical AND; both conditions must be true" },
    @{ Row = 16; Status = "Done"; DateDone = 46079; Language = "Python"; Time = "4 mins"; Difficulty = "Easy"; Notes = "Implemented flat slab calculation; used range comparisons; handled negative input; structured increasing conditions" }
)

foreach ($rd in $rowsData) {
    $r = $rd.Row

    # Column D: Status
    $ws.Range("D$r").Value = $rd.Status

    # Column E: Date Done — copy the date-formatted style from row 11 first
    # so the cell keeps the existing short-date number format, then set
    # the value itself (as the underlying date serial number).
    $ws.Range("E11").Copy($ws.Range("E$r"))
    $ws.Range("E$r").Value = $rd.DateDone

    # Columns F-I: Language, Time Taken, Difficulty, Notes
    $ws.Range("F$r").Value = $rd.Language
    $ws.Range("G$r").Value = $rd.Time
    $ws.Range("H$r").Value = $rd.Difficulty
    $ws.Range("I$r").Value = $rd.Notes
}

# Leave the cursor where the author last left it.
$ws.Range("H17").Select()

$wb.Save()
